$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows to append after the last existing data row (row 283).
$newDates = @("06-11-2021", "07-11-2021", "08-11-2021")

$startRow = 284
for ($i = 0; $i -lt $newDates.Length; $i++) {
    $r = $startRow + $i
    $cellA = $ws.Cells.Item($r, 1)
    $cellA.Value = "'" + $newDates[$i]
    $ws.Cells.Item($r, 2).Value = 449
    $ws.Cells.Item($r, 3).Value = 0
}
